$wb = $excel.ActiveWorkbook

$wsCosts = $wb.Worksheets.Item("costs")
$wsUtilities = $wb.Worksheets.Item("utilities")

# Insert a new "type" column (B) into both the costs and utilities sheets,
# shifting the existing cost/utility + variance columns one to the right.
$wsCosts.Columns("B").Insert()
$wsUtilities.Columns("B").Insert()

# --- costs sheet (was: state | cost | cost_variance) ---
# now:            state | type | cost | cost_variance
$wsCosts.Range("B1").Value = "type"
$wsCosts.Range("B2").Value = "static"
$wsCosts.Range("B3").Value = "static"
$wsCosts.Range("B4").Value = "static"
$wsCosts.Columns("B").ColumnWidth = $wsCosts.Columns("A").ColumnWidth

# --- utilities sheet (was: state | utility | utility_variance) ---
# now:                 state | type | utility | utility_variance
$wsUtilities.Range("B1").Value = "type"
$wsUtilities.Range("B2").Value = "static"
$wsUtilities.Range("B3").Value = "static"
$wsUtilities.Range("B4").Value = "static"
$wsUtilities.Columns("B").ColumnWidth = $wsUtilities.Columns("A").ColumnWidth

# Update selections to match the new layout, then make "costs" the active sheet/tab
# (selecting utilities first so the final .Select() below is what ends up active).
$wsUtilities.Range("B2:B4").Select()
$wsCosts.Range("B2:B4").Select()
